$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (A2 index 0)
$ws.Range("B2").Value = 17.28349325868628
$ws.Range("C2").Value = 14.017993757902
$ws.Range("E2").Value = 18.30425346541052
$ws.Range("F2").Value = 45.7800782152017
$ws.Range("G2").Value = 3.657627111594185
$ws.Range("I2").Value = 24.64803071066424
$ws.Range("J2").Value = 8.899799724719102

# Row 3 (A3 index 1)
$ws.Range("B3").Value = 16.68866467725297
$ws.Range("C3").Value = 13.39398741640595
$ws.Range("E3").Value = 18.20032058424491
$ws.Range("F3").Value = 45.41123953789438
$ws.Range("G3").Value = 3.662031050193564
$ws.Range("I3").Value = 24.65729904105183
$ws.Range("J3").Value = 8.924363253756715

# Row 4 (A4 index 2)
$ws.Range("B4").Value = 16.31740250672975
$ws.Range("C4").Value = 13.00015566069434
$ws.Range("E4").Value = 18.13943525298791
$ws.Range("F4").Value = 45.19867469584737
$ws.Range("G4").Value = 3.664869483473284
$ws.Range("I4").Value = 24.67153143813629
$ws.Range("J4").Value = 8.941251572216856

# Row 5 (A5 index 3)
$ws.Range("B5").Value = 16.16485327926443
$ws.Range("C5").Value = 12.83723901527339
$ws.Range("E5").Value = 18.11538187506929
$ws.Range("F5").Value = 45.11561718094524
$ws.Range("G5").Value = 3.666060109582046
$ws.Range("I5").Value = 24.67945972300969
$ws.Range("J5").Value = 8.948585769314764

# Row 6 (A6 index 4)
$ws.Range("B6").Value = 16.1394547129189
$ws.Range("C6").Value = 12.81004852503425
$ws.Range("E6").Value = 18.11143420392848
$ws.Range("F6").Value = 45.10204266121444
$ws.Range("G6").Value = 3.666259866407092
$ws.Range("I6").Value = 24.68090416411643
$ws.Range("J6").Value = 8.949830848675637

# Row 7 (A7 index 5)
$ws.Range("B7").Value = 16.31534990361521
$ws.Range("C7").Value = 12.99796796743485
$ws.Range("E7").Value = 18.13910776480383
$ws.Range("F7").Value = 45.19754003675131
$ws.Range("G7").Value = 3.664885403050554
$ws.Range("I7").Value = 24.67162977043658
$ws.Range("J7").Value = 8.941348656149966

# Row 8 (A8 index 6)
$ws.Range("B8").Value = 17.0797987580323
$ws.Range("C8").Value = 13.80520133210486
$ws.Range("E8").Value = 18.26782227716487
$ws.Range("F8").Value = 45.65006024952443
$ws.Range("G8").Value = 3.659117783812929
$ws.Range("I8").Value = 24.64944245610975
$ws.Range("J8").Value = 8.907893084909201

# Row 9 (A9 index 7)
$ws.Range("B9").Value = 18.52063251527898
$ws.Range("C9").Value = 15.29353185452762
$ws.Range("E9").Value = 18.54250883497081
$ws.Range("F9").Value = 46.64440638276633
$ws.Range("G9").Value = 3.648867190633973
$ws.Range("I9").Value = 24.67449490160437
$ws.Range("J9").Value = 8.856712241481855

# Row 10 (A10 index 8)
$ws.Range("B10").Value = 19.53136087307048
$ws.Range("C10").Value = 16.31821116967722
$ws.Range("E10").Value = 18.75655165777082
$ws.Range("F10").Value = 47.43517855923729
$ws.Range("G10").Value = 3.641972701159006
$ws.Range("I10").Value = 24.73565622397095
$ws.Range("J10").Value = 8.828031358049238

# Row 11 (A11 index 9)
$ws.Range("B11").Value = 19.97859465844454
$ws.Range("C11").Value = 16.76760965482848
$ws.Range("E11").Value = 18.8562798843855
$ws.Range("F11").Value = 47.80681832909287
$ws.Range("G11").Value = 3.63897243695019
$ws.Range("I11").Value = 24.77292912528343
$ws.Range("J11").Value = 8.816947713869929

# Row 12 (A12 index 10)
$ws.Range("B12").Value = 20.1459908006289
$ws.Range("C12").Value = 16.93525762141878
$ws.Range("E12").Value = 18.8943574567203
$ws.Range("F12").Value = 47.94915580263211
$ws.Range("G12").Value = 3.63785572564716
$ws.Range("I12").Value = 24.78841260955502
$ws.Range("J12").Value = 8.81303510847378

# Row 13 (A13 index 11)
$ws.Range("B13").Value = 20.11002858534645
$ws.Range("C13").Value = 16.8992658885665
$ws.Range("E13").Value = 18.88614329140362
$ws.Range("F13").Value = 47.91843122207066
$ws.Range("G13").Value = 3.638095367821272
$ws.Range("I13").Value = 24.78501693697656
$ws.Range("J13").Value = 8.813865069027665

# Row 14 (A14 index 12)
$ws.Range("B14").Value = 19.99240658408568
$ws.Range("C14").Value = 16.78145348549541
$ws.Range("E14").Value = 18.8594064364769
$ws.Range("F14").Value = 47.81849696925438
$ws.Range("G14").Value = 3.638880175983005
$ws.Range("I14").Value = 24.774175472767
$ws.Range("J14").Value = 8.816620106652632

# Row 15 (A15 index 13)
$ws.Range("B15").Value = 19.92009992473458
$ws.Range("C15").Value = 16.70895724842795
$ws.Range("E15").Value = 18.84306926347801
$ws.Range("F15").Value = 47.7574901932375
$ws.Range("G15").Value = 3.639363418912003
$ws.Range("I15").Value = 24.76771331695026
$ws.Range("J15").Value = 8.818344765321759

# Row 16 (A16 index 14)
$ws.Range("B16").Value = 19.5018666988559
$ws.Range("C16").Value = 16.28849421564623
$ws.Range("E16").Value = 18.75007939077895
$ws.Range("F16").Value = 47.41112164073854
$ws.Range("G16").Value = 3.64217150135671
$ws.Range("I16").Value = 24.73341142804979
$ws.Range("J16").Value = 8.828795397714702

# Row 17 (A17 index 15)
$ws.Range("B17").Value = 19.24196188379533
$ws.Range("C17").Value = 16.02617388749106
$ws.Range("E17").Value = 18.69361854220432
$ws.Range("F17").Value = 47.20161170297249
$ws.Range("G17").Value = 3.643928918545281
$ws.Range("I17").Value = 24.7147966331429
$ws.Range("J17").Value = 8.835710975444496

# Row 18 (A18 index 16)
$ws.Range("B18").Value = 19.09129821306427
$ws.Range("C18").Value = 15.87372517486202
$ws.Range("E18").Value = 18.66136741958339
$ws.Range("F18").Value = 47.08223518388468
$ws.Range("G18").Value = 3.644952554913134
$ws.Range("I18").Value = 24.70497893130566
$ws.Range("J18").Value = 8.839873311914742

# Row 19 (A19 index 17)
$ws.Range("B19").Value = 19.04008965245382
$ws.Range("C19").Value = 15.82184326632943
$ws.Range("E19").Value = 18.65048696181403
$ws.Range("F19").Value = 47.04201330942877
$ws.Range("G19").Value = 3.64530134621352
$ws.Range("I19").Value = 24.7018071766225
$ws.Range("J19").Value = 8.841314258933572

# Row 20 (A20 index 18)
$ws.Range("B20").Value = 19.26975182131355
$ws.Range("C20").Value = 16.05426167505334
$ws.Range("E20").Value = 18.69960593890458
$ws.Range("F20").Value = 47.22379841229207
$ws.Range("G20").Value = 3.643740513068301
$ws.Range("I20").Value = 24.71668611818755
$ws.Range("J20").Value = 8.834955670820897

# Row 21 (A21 index 19)
$ws.Range("B21").Value = 20.02700941015394
$ws.Range("C21").Value = 16.81612737954116
$ws.Range("E21").Value = 18.86725142691443
$ws.Range("F21").Value = 47.84780738158707
$ws.Range("G21").Value = 3.63864913279427
$ws.Range("I21").Value = 24.77732264800569
$ws.Range("J21").Value = 8.815803145986813

# Row 22 (A22 index 20)
$ws.Range("B22").Value = 20.5104275975112
$ws.Range("C22").Value = 17.29926294280796
$ws.Range("E22").Value = 18.97862730957752
$ws.Range("F22").Value = 48.26493963343854
$ws.Range("G22").Value = 3.635434770619967
$ws.Range("I22").Value = 24.82493516561482
$ws.Range("J22").Value = 8.804945311501113

# Row 23 (A23 index 21)
$ws.Range("B23").Value = 20.25351686671182
$ws.Range("C23").Value = 17.04279327742449
$ws.Range("E23").Value = 18.91902716331636
$ws.Range("F23").Value = 48.04149314835995
$ws.Range("G23").Value = 3.637140030501068
$ws.Range("I23").Value = 24.79879022720545
$ws.Range("J23").Value = 8.810587788880966

# Row 24 (A24 index 22)
$ws.Range("B24").Value = 19.25719184369551
$ws.Range("C24").Value = 16.04156827644568
$ws.Range("E24").Value = 18.69689838394675
$ws.Range("F24").Value = 47.21376444906784
$ws.Range("G24").Value = 3.643825649828806
$ws.Range("I24").Value = 24.71582912872979
$ws.Range("J24").Value = 8.835296563368185

# Row 25 (A25 index 23)
$ws.Range("B25").Value = 18.13844279847545
$ws.Range("C25").Value = 14.90228427140158
$ws.Range("E25").Value = 18.4659504362737
$ws.Range("F25").Value = 46.36444841427229
$ws.Range("G25").Value = 3.651527771837502
$ws.Range("I25").Value = 24.66027119844235
$ws.Range("J25").Value = 8.869000381205904
